$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Gratisflasche Amarone"
$ws.Range("O2").Value = "2022-09-14 21:00:20"
$ws.Range("A3").Value = "'6817200"
$ws.Range("B3").Value = "Zewa Wisch&amp;Weg decor 16 Rollen"
$ws.Range("C3").Value = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/zewa-wisch-weg-decor-16-rollen/p/6817200"
$ws.Range("D3").Value = "1ST"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = "Zewa"
$ws.Range("H3").Value = "'10.00"
$ws.Range("I3").Value = "10.00/1ST"
$ws.Range("J3").Value = "Preis pro 1 Stück"
$ws.Range("K3").Value = "'10.00"
$ws.Range("L3").Value = "1ST"
$ws.Range("M3").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
$ws.Range("N3").Value = "Zewa Wisch&amp;Weg decor 16 Rollen 54% Aktion 10.00 Schweizer Franken statt 22.00 Schweizer Franken"
$ws.Range("O3").Value = "2022-09-14 21:00:20"
$ws.Range("A4").Value = "'6283679"
$ws.Range("B4").Value = "Oecoplan Toilettenpapier Camomille weiss 4-lagig 6 Rollen"
$ws.Range("C4").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/oecoplan-toilettenpapier-camomille-weiss-4-lagig-6-rollen/p/6283679"
$ws.Range("D4").Value = "6Rol"
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 4
$ws.Range("H4").Value = "'4.50"
$ws.Range("I4").Value = "0.75/1Rol"
$ws.Range("J4").Value = "Preis pro 1 Rolle"
$ws.Range("K4").Value = "'0.75"
$ws.Range("L4").Value = "1Rol"
$ws.Range("M4").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N4").Value = "Oecoplan Toilettenpapier Camomille weiss 4-lagig 6 Rollen 4.50 Schweizer Franken"
$ws.Range("O4").Value = "2022-09-14 21:00:20"
$ws.Range("A5").Value = "'6695141"
$ws.Range("B5").Value = "Prix Garantie feuchtes Toilettenpapier 2x70 Stück"
$ws.Range("C5").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/prix-garantie-feuchtes-toilettenpapier-2x70-stueck/p/6695141"
$ws.Range("D5").Value = "140ST"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 3.5
$ws.Range("G5").Value = "Coop"
$ws.Range("H5").Value = "'2.50"
$ws.Range("I5").Value = "0.02/1ST"
$ws.Range("K5").Value = "'0.02"
$ws.Range("M5").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Range("N5").Value = "Prix Garantie feuchtes Toilettenpapier 2x70 Stück 2.50 Schweizer Franken"
$ws.Range("O5").Value = "2022-09-14 21:00:20"
$ws.Range("O6").Value = "2022-09-14 21:00:20"
$ws.Range("O7").Value = "2022-09-14 21:00:20"
$ws.Range("A8").Value = "'6834305"
$ws.Range("B8").Value = "Zewa Wisch&amp;Weg Haushaltspapier weiss 4 Rollen"
$ws.Range("C8").Value = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/zewa-wisch-weg-haushaltspapier-weiss-4-rollen/p/6834305"
$ws.Range("D8").Value = "192BLT"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = "Zewa"
$ws.Range("H8").Value = "'5.50"
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = ""
$ws.Range("M8").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
$ws.Range("N8").Value = "Zewa Wisch&amp;Weg Haushaltspapier weiss 4 Rollen 5.50 Schweizer Franken"
$ws.Range("O8").Value = "2022-09-14 21:00:20"
$ws.Range("A9").Value = "'4947421"
$ws.Range("B9").Value = "Oecoplan Taschentuch Calendula Box"
$ws.Range("C9").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/oecoplan-taschentuch-calendula-box/p/4947421"
$ws.Range("D9").Value = "80ST"
$ws.Range("E9").Value = 17
$ws.Range("F9").Value = 4
$ws.Range("G9").Value = "Coop"
$ws.Range("H9").Value = "'2.30"
$ws.Range("I9").Value = "0.03/1ST"
$ws.Range("J9").Value = "Preis pro 1 Stück"
$ws.Range("K9").Value = "'0.03"
$ws.Range("L9").Value = "1ST"
$ws.Range("M9").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Range("N9").Value = "Oecoplan Taschentuch Calendula Box 2.30 Schweizer Franken"
$ws.Range("O9").Value = "2022-09-14 21:00:20"
$ws.Range("A10").Value = "'6691348"
$ws.Range("B10").Value = "Super Soft Aloe Vera feucht FSC 4x  60ST"
$ws.Range("C10").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/super-soft-aloe-vera-feucht-fsc/p/6691348"
$ws.Range("D10").Value = "4x 60ST"
$ws.Range("E10").Value = 17
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = "Super Soft"
$ws.Range("H10").Value = "'7.65"
$ws.Range("I10").Value = "0.03/1ST"
$ws.Range("K10").Value = "'0.03"
$ws.Range("N10").Value = "Super Soft Aloe Vera feucht FSC 4x  60ST 35% Aktion 7.65 Schweizer Franken statt 11.80 Schweizer Franken"
$ws.Range("O10").Value = "2022-09-14 21:00:20"
$ws.Range("A11").Value = "'6724076"
$ws.Range("B11").Value = "Oecoplan feuchtes Toilettenpapier Duckies natural 40 Stück"
$ws.Range("C11").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/oecoplan-feuchtes-toilettenpapier-duckies-natural-40-stueck/p/6724076"
$ws.Range("D11").Value = "40ST"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 4.5
$ws.Range("G11").Value = "Duckies"
$ws.Range("H11").Value = "'2.95"
$ws.Range("I11").Value = "0.07/1ST"
$ws.Range("J11").Value = "Preis pro 1 Stück"
$ws.Range("K11").Value = "'0.07"
$ws.Range("L11").Value = "1ST"
$ws.Range("M11").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Range("N11").Value = "Oecoplan feuchtes Toilettenpapier Duckies natural 40 Stück 2.95 Schweizer Franken"
$ws.Range("O11").Value = "2022-09-14 21:00:20"
$ws.Range("A12").Value = "'6498157"
$ws.Range("B12").Value = "subito Haushaltspapier weiss 2 Rollen"
$ws.Range("C12").Value = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/subito-haushaltspapier-weiss-2-rollen/p/6498157"
$ws.Range("D12").Value = "100BLT"
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = "subito"
$ws.Range("H12").Value = "'2.75"
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
$ws.Range("N12").Value = "subito Haushaltspapier weiss 2 Rollen 2.75 Schweizer Franken"
$ws.Range("O12").Value = "2022-09-14 21:00:20"
$ws.Range("O13").Value = "2022-09-14 21:00:20"
$ws.Range("O14").Value = "2022-09-14 21:00:20"
$ws.Range("O15").Value = "2022-09-14 21:00:20"
$ws.Range("O16").Value = "2022-09-14 21:00:20"
$ws.Range("O17").Value = "2022-09-14 21:00:20"
$ws.Range("O18").Value = "2022-09-14 21:00:20"
$ws.Range("O19").Value = "2022-09-14 21:00:20"
$ws.Range("O20").Value = "2022-09-14 21:00:20"
$ws.Range("O21").Value = "2022-09-14 21:00:20"
$ws.Range("O22").Value = "2022-09-14 21:00:20"
$ws.Range("O23").Value = "2022-09-14 21:00:20"
$ws.Range("A24").Value = "'6866582"
$ws.Range("B24").Value = "Tempo Bamboo Eco Box"
$ws.Range("C24").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/tempo-bamboo-eco-box/p/6866582"
$ws.Range("D24").Value = "90ST"
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 5
$ws.Range("G24").Value = "Tempo"
$ws.Range("H24").Value = "'3.95"
$ws.Range("I24").Value = "0.04/1ST"
$ws.Range("K24").Value = "'0.04"
$ws.Range("M24").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Range("N24").Value = "Tempo Bamboo Eco Box 3.95 Schweizer Franken"
$ws.Range("O24").Value = "2022-09-14 21:00:20"
$ws.Range("A25").Value = "'4687972"
$ws.Range("B25").Value = "Taschentücher Strong 10x10 Stück"
$ws.Range("C25").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/taschentuecher-strong-10x10-stueck/p/4687972"
$ws.Range("D25").Value = "10ST"
$ws.Range("G25").Value = "Super Silk"
$ws.Range("H25").Value = "'1.80"
$ws.Range("I25").Value = "0.18/1ST"
$ws.Range("K25").Value = "'0.18"
$ws.Range("N25").Value = "Taschentücher Strong 10x10 Stück 1.80 Schweizer Franken"
$ws.Range("O25").Value = "2022-09-14 21:00:20"
$ws.Range("A26").Value = "'6996030"
$ws.Range("B26").Value = "Tela Viva Haushaltspapier 3-lagig 4 Rollen"
$ws.Range("C26").Value = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/tela-viva-haushaltspapier-3-lagig-4-rollen/p/6996030"
$ws.Range("D26").Value = "200BLT"
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = "Tela"
$ws.Range("H26").Value = "'5.95"
$ws.Range("I26").Value = ""
$ws.Range("J26").Value = ""
$ws.Range("K26").Value = ""
$ws.Range("L26").Value = ""
$ws.Range("M26").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
$ws.Range("N26").Value = "Tela Viva Haushaltspapier 3-lagig 4 Rollen 5.95 Schweizer Franken"
$ws.Range("O26").Value = "2022-09-14 21:00:20"
$ws.Range("A27").Value = "'3874909"
$ws.Range("B27").Value = "Oecoplan Papiertaschentücher Special-Edition Calendula 30x10 Stück"
$ws.Range("C27").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/oecoplan-papiertaschentuecher-special-edition-calendula-30x10-stueck/p/3874909"
$ws.Range("D27").Value = "30ST"
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = "Coop"
$ws.Range("H27").Value = "'3.65"
$ws.Range("I27").Value = "0.12/1ST"
$ws.Range("J27").Value = "Preis pro 1 Stück"
$ws.Range("K27").Value = "'0.12"
$ws.Range("L27").Value = "1ST"
$ws.Range("M27").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Range("N27").Value = "Oecoplan Papiertaschentücher Special-Edition Calendula 30x10 Stück 20% Aktion 3.65 Schweizer Franken statt 4.60 Schweizer Franken"
$ws.Range("O27").Value = "2022-09-14 21:00:20"
$ws.Range("A28").Value = "'6800946"
$ws.Range("B28").Value = "Hipp Natural zart duftend 3x48 Stück"
$ws.Range("C28").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/hipp-natural-zart-duftend-3x48-stueck/p/6800946"
$ws.Range("D28").Value = "144ST"
$ws.Range("E28").Value = 2
$ws.Range("G28").Value = "Hipp"
$ws.Range("H28").Value = "'6.65"
$ws.Range("I28").Value = "0.05/1ST"
$ws.Range("K28").Value = "'0.05"
$ws.Range("M28").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Range("N28").Value = "Hipp Natural zart duftend 3x48 Stück 33% Aktion 6.65 Schweizer Franken statt 9.95 Schweizer Franken"
$ws.Range("O28").Value = "2022-09-14 21:00:20"
$ws.Range("A29").Value = "'6996129"
$ws.Range("B29").Value = "Tela Toilettenpapier Futura 3-lagig 9 Rollen"
$ws.Range("C29").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/tela-toilettenpapier-futura-3-lagig-9-rollen/p/6996129"
$ws.Range("D29").Value = "9Rol"
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = "Tela"
$ws.Range("H29").Value = "'8.40"
$ws.Range("I29").Value = "0.93/1Rol"
$ws.Range("J29").Value = "Preis pro 1 Rolle"
$ws.Range("K29").Value = "'0.93"
$ws.Range("L29").Value = "1Rol"
$ws.Range("M29").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N29").Value = "Tela Toilettenpapier Futura 3-lagig 9 Rollen 8.40 Schweizer Franken"
$ws.Range("O29").Value = "2022-09-14 21:00:20"
$ws.Range("A30").Value = "'6868354"
$ws.Range("B30").Value = "Tempo Bamboo Eco"
$ws.Range("C30").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/tempo-bamboo-eco/p/6868354"
$ws.Range("D30").Value = "12ST"
$ws.Range("E30").Value = 1
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = "Tempo"
$ws.Range("H30").Value = "'3.95"
$ws.Range("I30").Value = "0.33/1ST"
$ws.Range("J30").Value = "Preis pro 1 Stück"
$ws.Range("K30").Value = "'0.33"
$ws.Range("L30").Value = "1ST"
$ws.Range("M30").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Range("N30").Value = "Tempo Bamboo Eco 3.95 Schweizer Franken"
$ws.Range("O30").Value = "2022-09-14 21:00:20"
$ws.Range("A31").Value = "'6283677"
$ws.Range("B31").Value = "Oecoplan Goldmelisse blau 3-lagig 32 Rollen"
$ws.Range("C31").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/oecoplan-goldmelisse-blau-3-lagig-32-rollen/p/6283677"
$ws.Range("D31").Value = "32Rol"
$ws.Range("F31").Value = 5
$ws.Range("G31").Value = "Coop"
$ws.Range("H31").Value = "'14.80"
$ws.Range("I31").Value = "0.46/1Rol"
$ws.Range("J31").Value = "Preis pro 1 Rolle"
$ws.Range("K31").Value = "'0.46"
$ws.Range("L31").Value = "1Rol"
$ws.Range("M31").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N31").Value = "Oecoplan Goldmelisse blau 3-lagig 32 Rollen 30% Aktion 14.80 Schweizer Franken statt 21.20 Schweizer Franken"
$ws.Range("O31").Value = "2022-09-14 21:00:20"
$ws.Range("O32").Value = "2022-09-14 21:00:20"
$ws.Range("A33").Value = "'6636712"
$ws.Range("B33").Value = "Pampers Coconut Pure 42 Feuchttücher"
$ws.Range("C33").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/pampers-coconut-pure-42-feuchttuecher/p/6636712"
$ws.Range("D33").Value = "42ST"
$ws.Range("E33").Value = 2
$ws.Range("F33").Value = 3.5
$ws.Range("G33").Value = "Pampers"
$ws.Range("H33").Value = "'4.95"
$ws.Range("I33").Value = "0.12/1ST"
$ws.Range("J33").Value = "Preis pro 1 Stück"
$ws.Range("K33").Value = "'0.12"
$ws.Range("L33").Value = "1ST"
$ws.Range("M33").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Range("N33").Value = "Pampers Coconut Pure 42 Feuchttücher 4.95 Schweizer Franken"
$ws.Range("O33").Value = "2022-09-14 21:00:20"
$ws.Range("A34").Value = "'6727164"
$ws.Range("B34").Value = "Wetties Allzwecktücher 80Stück"
$ws.Range("C34").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/wetties-allzwecktuecher-80stueck/p/6727164"
$ws.Range("D34").Value = "80ST"
$ws.Range("E34").Value = 1
$ws.Range("F34").Value = 5
$ws.Range("G34").Value = "Coop"
$ws.Range("H34").Value = "'3.95"
$ws.Range("I34").Value = "0.05/1ST"
$ws.Range("K34").Value = "'0.05"
$ws.Range("N34").Value = "Wetties Allzwecktücher 80Stück 3.95 Schweizer Franken"
$ws.Range("O34").Value = "2022-09-14 21:00:20"

$ws.Rows.Item(35).Delete()